$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# win / total win counters: 4 new rows (25-28) mirroring the existing
# "message box" rows (21-24) layout: label in column A, numeric
# width/height/x/y in columns B:E.
# ---------------------------------------------------------------------

# -- values ------------------------------------------------------------
$ws.Range("A25").Value = "win image"
$ws.Range("B25").Value = 80
$ws.Range("C25").Value = 35
$ws.Range("D25").Value = 172
$ws.Range("E25").Value = 97

$ws.Range("A26").Value = "total win image"
$ws.Range("B26").Value = 80
$ws.Range("C26").Value = 35
$ws.Range("D26").Value = 330
$ws.Range("E26").Value = 97

$ws.Range("A27").Value = "win label"
$ws.Range("D27").Value = 270
$ws.Range("E27").Value = 110

$ws.Range("A28").Value = "total win label"
$ws.Range("D28").Value = 420
$ws.Range("E28").Value = 110

# -- formatting ----------------------------------------------------------
# Column A labels: copy the look of the existing label cells (A21:A24)
$ws.Range("A21:A24").Copy()
$ws.Range("A25:A28").PasteSpecial(-4122)

# Columns B:E on rows 25/26: copy the centered data-cell look used by
# the rows above (B21:E21)
$ws.Range("B21:E21").Copy()
$ws.Range("B25:E26").PasteSpecial(-4122)

# Columns B:E on rows 27/28: same bordered/filled look, but centered
$ws.Range("B21:E21").Copy()
$ws.Range("B27:E28").PasteSpecial(-4122)
$ws.Range("B27:E28").HorizontalAlignment = -4108

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# sheet view / selection bookkeeping
# ---------------------------------------------------------------------
$ws.Range("A1").CurrentRegion | Out-Null
$ws.Activate()
$ws.Range("B27").Select()
$av = $excel.ActiveWindow
$av.ScrollRow = 16
